# Trading update: 2026-02-18 11:06:16
# Two new trades landed:
#   Trade #11 - MarketMaking - OPEN
#   Trade #12 - momentum     - OPEN
# A new "momentum" strategy sheet is created (positioned right after
# "All Trades" and before "MarketMaking"), the "MarketMaking" strategy
# sheet is refreshed to show only its newest trade, and the master
# "All Trades" log gets two new rows appended (with the previous two
# "latest trade" rows' extra bookkeeping columns cleared out, since
# they are no longer the newest trade for their strategy).

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

$header = @("Trade #","Date","Time","Strategy","Side","Entry Price","Exit Price","Status","P&L %","P&L $","Capital After","Entry Slippage (bps)","Exit Slippage (bps)","Confidence","Entry Reason","Exit Reason","Duration (min)")

# ---------------------------------------------------------------
# 1) Recreate "MarketMaking" (dropped, then re-added after the new
#    "momentum" sheet) so the internal sheet id ordering mirrors a
#    brand-new strategy tab being inserted ahead of it, and build
#    the new "momentum" sheet positioned between "All Trades" and
#    "MarketMaking".
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$oldMarketMaking = $wb.Worksheets.Item("MarketMaking")
[void]$oldMarketMaking.Delete()

$momentum = $wb.Worksheets.Add($null, $allTrades)
$momentum.Name = "momentum"
$momentum.PageSetup.LeftMargin = 0.75 * 72
$momentum.PageSetup.RightMargin = 0.75 * 72
$momentum.PageSetup.TopMargin = 1 * 72
$momentum.PageSetup.BottomMargin = 1 * 72
$momentum.PageSetup.HeaderMargin = 0.5 * 72
$momentum.PageSetup.FooterMargin = 0.5 * 72

$marketMaking = $wb.Worksheets.Add($null, $momentum)
$marketMaking.Name = "MarketMaking"
$marketMaking.PageSetup.LeftMargin = 0.75 * 72
$marketMaking.PageSetup.RightMargin = 0.75 * 72
$marketMaking.PageSetup.TopMargin = 1 * 72
$marketMaking.PageSetup.BottomMargin = 1 * 72
$marketMaking.PageSetup.HeaderMargin = 0.5 * 72
$marketMaking.PageSetup.FooterMargin = 0.5 * 72

# ---------------------------------------------------------------
# 2) Populate the "momentum" sheet with its header + newest trade.
# ---------------------------------------------------------------
for ($i = 0; $i -lt $header.Length; $i++) {
    $momentum.Cells.Item(1, $i + 1).Value = $header[$i]
}

$momentum.Cells.Item(2, 1).Value = 12
$momentum.Cells.Item(2, 2).Value = "'2026-02-18"
$momentum.Cells.Item(2, 3).Value = "11:05:20"
$momentum.Cells.Item(2, 4).Value = "momentum"
$momentum.Cells.Item(2, 5).Value = "UP"
$momentum.Cells.Item(2, 6).Value = 0.49
$momentum.Cells.Item(2, 8).Value = "OPEN"
$momentum.Cells.Item(2, 9).Value = 0
$momentum.Cells.Item(2, 10).Value = 0
$momentum.Cells.Item(2, 11).Value = 100
$momentum.Cells.Item(2, 12).Value = 0
$momentum.Cells.Item(2, 13).Value = 0
$momentum.Cells.Item(2, 14).Value = 0.9
$momentum.Cells.Item(2, 15).Value = "Upward momentum: 22.222% over 5 samples"
$momentum.Cells.Item(2, 17).Value = 0

# ---------------------------------------------------------------
# 3) Populate the refreshed "MarketMaking" sheet with its header +
#    newest trade (#11).
# ---------------------------------------------------------------
for ($i = 0; $i -lt $header.Length; $i++) {
    $marketMaking.Cells.Item(1, $i + 1).Value = $header[$i]
}

$marketMaking.Cells.Item(2, 1).Value = 11
$marketMaking.Cells.Item(2, 2).Value = "'2026-02-18"
$marketMaking.Cells.Item(2, 3).Value = "11:05:05"
$marketMaking.Cells.Item(2, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(2, 5).Value = "UP"
$marketMaking.Cells.Item(2, 6).Value = 0.46
$marketMaking.Cells.Item(2, 8).Value = "OPEN"
$marketMaking.Cells.Item(2, 9).Value = 0
$marketMaking.Cells.Item(2, 10).Value = 0
$marketMaking.Cells.Item(2, 11).Value = 100
$marketMaking.Cells.Item(2, 12).Value = 0
$marketMaking.Cells.Item(2, 13).Value = 0
$marketMaking.Cells.Item(2, 14).Value = 0.6
$marketMaking.Cells.Item(2, 15).Value = "Normal spread capture: 202 bps"
$marketMaking.Cells.Item(2, 17).Value = 0

# ---------------------------------------------------------------
# 4) Update the "All Trades" master log: the previous "latest
#    trade" bookkeeping columns (K:O, Q) on rows 10 and 11 are no
#    longer the newest for their strategy, so clear them out (Exit
#    Price column G flips from blank to 0, matching the other
#    superseded rows), then append the two new trades.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("All Trades")

$ws.Range("G10").Value = 0
$ws.Range("K10:O10").ClearContents()
$ws.Range("Q10").ClearContents()

$ws.Range("G11").Value = 0
$ws.Range("K11:O11").ClearContents()
$ws.Range("Q11").ClearContents()

# New row 12: Trade #11 (MarketMaking)
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "'2026-02-18"
$ws.Cells.Item(12, 3).Value = "11:05:05"
$ws.Cells.Item(12, 4).Value = "MarketMaking"
$ws.Cells.Item(12, 5).Value = "UP"
$ws.Cells.Item(12, 6).Value = 0.46
$ws.Cells.Item(12, 8).Value = "OPEN"
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 100
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = 0.6
$ws.Cells.Item(12, 15).Value = "Normal spread capture: 202 bps"
$ws.Cells.Item(12, 17).Value = 0

# New row 13: Trade #12 (momentum)
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "'2026-02-18"
$ws.Cells.Item(13, 3).Value = "11:05:20"
$ws.Cells.Item(13, 4).Value = "momentum"
$ws.Cells.Item(13, 5).Value = "UP"
$ws.Cells.Item(13, 6).Value = 0.49
$ws.Cells.Item(13, 8).Value = "OPEN"
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 100
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 0.9
$ws.Cells.Item(13, 15).Value = "Upward momentum: 22.222% over 5 samples"
$ws.Cells.Item(13, 17).Value = 0

# ---------------------------------------------------------------
# 5) Restore the originally active sheet/selection (creating the
#    new sheets shifts focus onto them as a side effect).
# ---------------------------------------------------------------
[void]$wb.Worksheets.Item("Summary").Activate()
[void]$wb.Worksheets.Item("Summary").Range("A1").Select()
